$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the fantasy roster table (A2:C16) with updated player/position/team
# assignments. Rows 17-19 (Giannis Antetokounmpo, Fred VanVleet, RJ Barrett)
# are unchanged and left as-is.
$ws.Range("A2").Value = "Anthony Edwards"
$ws.Range("B2").Value = "SG,SF"
$ws.Range("C2").Value = "Minnesota Timberwolves"

$ws.Range("A3").Value = "James Harden"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "LA Clippers"

$ws.Range("A4").Value = "Anfernee Simons"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Portland Trail Blazers"

$ws.Range("A5").Value = "Jayson Tatum"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "Boston Celtics"

$ws.Range("A6").Value = "Paul George"
$ws.Range("B6").Value = "SG,SF,PF"
$ws.Range("C6").Value = "Philadelphia 76ers"

$ws.Range("A7").Value = "Amen Thompson"
$ws.Range("B7").Value = "SG,SF,PF"
$ws.Range("C7").Value = "Houston Rockets"

$ws.Range("A8").Value = "Zion Williamson"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "New Orleans Pelicans"

$ws.Range("A9").Value = "Nicolas Claxton"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Brooklyn Nets"

$ws.Range("A10").Value = "Zach Edey"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Memphis Grizzlies"

$ws.Range("A11").Value = "Jaren Jackson Jr."
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Memphis Grizzlies"

$ws.Range("A12").Value = "Ivica Zubac"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "LA Clippers"

$ws.Range("A13").Value = "Bradley Beal"
$ws.Range("B13").Value = "PG,SG,SF"
$ws.Range("C13").Value = "Phoenix Suns"

$ws.Range("A14").Value = "Bobby Portis"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Milwaukee Bucks"

$ws.Range("A15").Value = "Keyonte George"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Utah Jazz"

$ws.Range("A16").Value = "Draymond Green"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Golden State Warriors"
